$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to match the source file name
$ws.Name = "study_background_traits-ALL"

# Fix the background trait value for GCST002887 (row 183, column E / "Background")
# Remove the extra trailing "||" split delimiter
$ws.Range("E183").Value = "unipolar depression|| schizophrenia|| bipolar disorder"

# Highlight the corrected row (A183:E183) with red font on white fill
$rng = $ws.Range("A183:E183")
$rng.Font.Color = 255
$rng.Interior.Color = 16777215

# Update the sheet view: scroll the frozen pane and move the selection
$ws.Application.ActiveWindow.ScrollRow = 170
$ws.Range("F183").Select()
